$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: price updated
$ws.Range("E2").Value = 1

# Row 3: Base3 changed from DEF789 to DEF345; ValidTo (D3) cleared
$ws.Range("B3").Value = "DEF345"
$ws.Range("D3").ClearContents()

# Row 4: now represents ABC123 / DEF789, with both ValidFrom and ValidTo, price 1
$ws.Range("A4").Value = "ABC123"
$ws.Range("B4").Value = "DEF789"
$ws.Range("C4").Value = "2021-01-01T00:00:00+00:00"
$ws.Range("D4").Value = "2023-01-01T00:00:00+00:00"
$ws.Range("E4").Value = 1

# Row 5 (new example row): DEF345 / QWE111, ValidTo only, price 1.5
$ws.Range("A5").Value = "DEF345"
$ws.Range("B5").Value = "QWE111"
$ws.Range("D5").Value = "2023-01-01T00:00:00+00:00"
$ws.Range("E5").Value = 1.5

$ws.Range("E4").Select() | Out-Null
